# Apply the StructureDefinition metadata refresh (5.0.0 -> 6.0.0 release):
#  - Metadata sheet: bump Version / Date, fill in Publisher, replace the
#    duplicated "Contact" row with a "Jurisdiction" row, and drop the
#    second (duplicate) Contact row entirely (table shrinks by one row).
#  - Elements sheet: give the root Extension row a real Short/Definition
#    instead of the generic placeholder text.

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refreshed publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicate "Contact" row -> becomes "Jurisdiction"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was the second duplicate "Contact" row; delete it entirely so
# every row below (old Description ... Context) shifts up by one and the
# table dimension shrinks from A1:B21 to A1:B20.
$meta.Rows("11").Delete()

# --- Elements sheet --------------------------------------------------------

# Root Extension row (row 2): Short / Definition now reflect the actual
# extension instead of the generic "Extension" / "An Extension" text.
$elements.Range("K2").Value = "SubscriberId Encrypted"
$elements.Range("L2").Value = "Encrypted identifier of the subscriber or family"
